{"js": "// Korean Welcome.docx translation-pass edit:\n// Trim the trailing clause \", until Smartcash reaches a considerable\n// market cap\" from the ASIC/mining paragraph so the sentence now ends\n// \"...for quite some time.\" (the rest of the paragraph, and its\n// formatting, is unchanged).\n\nconst body = context.document.body;\n\nconst results = body.search(\n  \", until Smartcash reaches a considerable market cap.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Replace the clause (including the leading comma) with just a\n  // period so the sentence reads \"...for quite some time.\"\n  results.items[i].insertText(\".\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Korean Welcome.docx translation-pass edit:\n# Trim the trailing clause \", until Smartcash reaches a considerable\n# market cap\" from the ASIC/mining paragraph so the sentence now ends\n# \"...for quite some time.\" (rest of paragraph/formatting unchanged).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \", until Smartcash reaches a considerable market cap.\"\n$find.Replacement.Text = \".\"\n$find.Execute(\n    $find.Text,     # FindText\n    $false,         # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2               # Replace (wdReplaceAll)\n)\n"}
